# Add a new Job Posting (Job_Id = JD_001) to Sheet1: a bold / boxed /
# centered header row (A1:K1) plus the corresponding data row (row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
$headers = @(
    "Job_Id",
    "Job_Title",
    "Job_Description",
    "Total_Years_Min_Exp",
    "Total_Years_Max_Exp",
    "Work_Mode",
    "Job_Location",
    "LinkedIn_Poster",
    "LinkedIn_Posted",
    "Resume_received",
    "Resume_downloaded"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Build the header look (bold font, thin box border, centered horizontally,
# top-aligned vertically) on an out-of-the-way scratch cell first, then copy
# just the formatting onto the header range in a single paste so the whole
# header gets one consistent cell style instead of layering several.
$scratch = $ws.Range("Z100")
$scratch.Font.Bold = $true
$scratch.Borders.LineStyle = 1
$scratch.HorizontalAlignment = -4108   # xlCenter
$scratch.VerticalAlignment = -4160     # xlTop

$headerRange = $ws.Range("A1:K1")
$scratch.Copy()
$headerRange.PasteSpecial(-4122)       # xlPasteFormats
$excel.CutCopyMode = $false
$scratch.Clear()

# --- Data row (row 2) -----------------------------------------------------
$ws.Range("A2").Value = "JD_001"
$ws.Range("B2").Value = "Junior RPA Developer"
$ws.Range("C2").Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = "Remote"
$ws.Range("G2").Value = "Bengaluru, Karnataka, India"
